# Applies the cryptos.xlsx price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to keep a literal text representation (e.g. "0.620")
    # instead of being auto-coerced into a number that would drop trailing
    # zeros / formatting, mirroring the inline string cells in the workbook.
    $range.NumberFormat = "@"
    $range.Value = $value
}

$ws.Range('D2').Value = '40.732.51'
$ws.Range('E2').Value = '  -7.07%  '
$ws.Range('D3').Value = '2.187.23'
$ws.Range('E3').Value = '  -7.20%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue $ws.Range('D5') '241.03'
$ws.Range('E5').Value = '  +0.49%  '
Set-TextValue $ws.Range('D6') '0.619'
$ws.Range('E6').Value = '  -7.70%  '
Set-TextValue $ws.Range('D7') '68.61'
$ws.Range('E7').Value = '  -7.76%  '
$ws.Range('E8').Value = '  +0.17%  '
Set-TextValue $ws.Range('D9') '0.535'
$ws.Range('E9').Value = '  -12.01%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D10') '0.0943'
$ws.Range('E10').Value = '  -7.68%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D11') '36.37'
$ws.Range('E11').Value = '  +1.80%  '
Set-TextValue $ws.Range('D12') '57.56'
$ws.Range('E12').Value = '  -5.52%  '
$ws.Range('E13').Value = '  -4.35%  '
Set-TextValue $ws.Range('D14') '6.56'
$ws.Range('E14').Value = '  -9.44%  '
$ws.Range('D15').Value = '2.509.64'
$ws.Range('E15').Value = '  -7.36%  '
Set-TextValue $ws.Range('D16') '14.59'
$ws.Range('E16').Value = '  -10.61%  '
Set-TextValue $ws.Range('D17') '0.828'
$ws.Range('E17').Value = '  -9.87%  '
$ws.Range('D18').Value = '2.187.43'
$ws.Range('E18').Value = '  -7.32%  '
$ws.Range('D19').Value = '40.665.93'
$ws.Range('E19').Value = '  -7.14%  '
$ws.Range('D20').Value = '0.0₃0936'
$ws.Range('E20').Value = '  -9.16%  '
Set-TextValue $ws.Range('D21') '72.27'
$ws.Range('E21').Value = '  -7.21%  '
Set-TextValue $ws.Range('D22') '6.03'
$ws.Range('E22').Value = '  -8.11%  '
Set-TextValue $ws.Range('D23') '229.18'
$ws.Range('E23').Value = '  -9.58%  '
$ws.Range('E24').Value = '  +7.33%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -5.14%  '
$ws.Range('E27').Value = '  -4.35%  '
Set-TextValue $ws.Range('D28') '2.17'
$ws.Range('E28').Value = '  -5.17%  '
Set-TextValue $ws.Range('D29') '9.65'
$ws.Range('E29').Value = '  -8.06%  '
Set-TextValue $ws.Range('D30') '168.21'
$ws.Range('E30').Value = '  -4.15%  '
Set-TextValue $ws.Range('D31') '20.16'
$ws.Range('E31').Value = '  -9.69%  '
$ws.Range('E32').Value = '  -10.01%  '
Set-TextValue $ws.Range('D33') '0.122'
$ws.Range('E33').Value = '  -8.20%  '
Set-TextValue $ws.Range('D34') '0.0696'
$ws.Range('E34').Value = '  -6.98%  '
Set-TextValue $ws.Range('D35') '5.10'
$ws.Range('E35').Value = '  -5.17%  '
Set-TextValue $ws.Range('D36') '4.54'
$ws.Range('E36').Value = '  -10.37%  '
Set-TextValue $ws.Range('D37') '3.82'
$ws.Range('E37').Value = '  +0.65%  '
Set-TextValue $ws.Range('D38') '23.26'
$ws.Range('E38').Value = '  +14.88%  '
$ws.Range('E39').Value = '  -6.86%  '
$ws.Range('E40').Value = '  -4.46%  '
Set-TextValue $ws.Range('D41') '5.77'
$ws.Range('E41').Value = '  -12.50%  '
Set-TextValue $ws.Range('D42') '62.76'
$ws.Range('E42').Value = '  -3.92%  '
Set-TextValue $ws.Range('D43') '4.82'
$ws.Range('E43').Value = '  -12.23%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D44') '8.55'
$ws.Range('E44').Value = '  -5.49%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D45') '0.190'
$ws.Range('E45').Value = '  -6.68%  '
$ws.Range('E46').Value = '  -0.18%  '
Set-TextValue $ws.Range('D47') '0.0977'
$ws.Range('E47').Value = '  -8.07%  '
$ws.Range('E48').Value = '  +2.48%  '
Set-TextValue $ws.Range('D49') '10.24'
$ws.Range('E49').Value = '  +6.19%  '
Set-TextValue $ws.Range('D50') '1.15'
$ws.Range('E50').Value = '  -6.64%  '
$ws.Range('E51').Value = '  -6.31%  '

Write-Host "Applied 95 cell updates"
